$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# For both the "zh-cn" and "de-de" worksheets:
#   * Update the Status column (B) for the two file rows to
#     "Handed back: in sync with en-US"
#   * Populate "Latest Target File" (E) and "Latest Handback File" (F) with
#     hyperlinks that mirror the existing Source File (A) / Latest Handoff
#     File (C) hyperlinks
#   * Record the new "Latest Handback DateTime" (G) for the two file rows
# ---------------------------------------------------------------------------

$statusHandedBack = "Handed back: in sync with en-US"

# Matches the workbook's existing custom "HyperLink" cell style (underline +
# RGB 0x6495ED font color) that is already applied to columns A/C.
function Format-AsExistingHyperlinkStyle($range) {
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

function Apply-HandbackReport {
    param(
        [string]$sheetName,
        [string]$mdUrl1,
        [string]$xlfUrl1,
        [string]$mdUrl2,
        [string]$xlfUrl2,
        [string]$handbackDateTime
    )

    $ws = $wb.Worksheets.Item($sheetName)

    # --- Row 2 (471bcea5-e21d-46bc-9a06-3774b23cce26.md) -------------------
    $ws.Range("B2").Value = $statusHandedBack

    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl1, [Type]::Missing, [Type]::Missing, "471bcea5-e21d-46bc-9a06-3774b23cce26.md")
    Format-AsExistingHyperlinkStyle $ws.Range("E2")

    $ws.Hyperlinks.Add($ws.Range("F2"), $xlfUrl1, [Type]::Missing, [Type]::Missing, "471bcea5-e21d-46bc-9a06-3774b23cce26.803d5d1e245744638b6d29189e46fe7612de3601.$sheetName.xlf")
    Format-AsExistingHyperlinkStyle $ws.Range("F2")

    $ws.Range("G2").Value = $handbackDateTime

    # --- Row 3 (de0ea274-0b1c-4ac4-be6e-532d1adf082a.md) --------------------
    $ws.Range("B3").Value = $statusHandedBack

    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl2, [Type]::Missing, [Type]::Missing, "de0ea274-0b1c-4ac4-be6e-532d1adf082a.md")
    Format-AsExistingHyperlinkStyle $ws.Range("E3")

    $ws.Hyperlinks.Add($ws.Range("F3"), $xlfUrl2, [Type]::Missing, [Type]::Missing, "de0ea274-0b1c-4ac4-be6e-532d1adf082a.04638db6c3a971fa0468bce48ee252a1981fee14.$sheetName.xlf")
    Format-AsExistingHyperlinkStyle $ws.Range("F3")

    $ws.Range("G3").Value = $handbackDateTime
}

# zh-cn sheet -----------------------------------------------------------
Apply-HandbackReport `
    "zh-cn" `
    "https://github.com/OpenLocalizationTest/oltest/blob/100916f69cefc51994f3f8dbb38ab33e33428571/e2e/471bcea5-e21d-46bc-9a06-3774b23cce26.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/390c0c617562bd8f3cace074b8b22ccb2e3ba5a7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/471bcea5-e21d-46bc-9a06-3774b23cce26.803d5d1e245744638b6d29189e46fe7612de3601.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/100916f69cefc51994f3f8dbb38ab33e33428571/e2e/de0ea274-0b1c-4ac4-be6e-532d1adf082a.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/390c0c617562bd8f3cace074b8b22ccb2e3ba5a7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/de0ea274-0b1c-4ac4-be6e-532d1adf082a.04638db6c3a971fa0468bce48ee252a1981fee14.zh-cn.xlf" `
    "2016-02-24 08:57:20"

# de-de sheet -----------------------------------------------------------
Apply-HandbackReport `
    "de-de" `
    "https://github.com/OpenLocalizationTest/oltest/blob/100916f69cefc51994f3f8dbb38ab33e33428571/e2e/471bcea5-e21d-46bc-9a06-3774b23cce26.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/58ae4e6607c1723d177a02c4c34177a58f5a91fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/471bcea5-e21d-46bc-9a06-3774b23cce26.803d5d1e245744638b6d29189e46fe7612de3601.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/100916f69cefc51994f3f8dbb38ab33e33428571/e2e/de0ea274-0b1c-4ac4-be6e-532d1adf082a.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/58ae4e6607c1723d177a02c4c34177a58f5a91fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/de0ea274-0b1c-4ac4-be6e-532d1adf082a.04638db6c3a971fa0468bce48ee252a1981fee14.de-de.xlf" `
    "2016-02-24 08:57:42"

$wb.Save()
